$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three test/quiz column headers
$ws.Range("B1").Value = "Notes devoir-test octobre"
$ws.Range("C1").Value = "Notes devoir-test chapitre 4"
$ws.Range("D1").Value = "Notes d'examen trismestriel"

# Widen the columns to fit the new, longer header text
$ws.Columns.Item(1).ColumnWidth = 30.6667
$ws.Columns.Item(2).ColumnWidth = 27.6667
$ws.Columns.Item(3).ColumnWidth = 30.6667
$ws.Columns.Item(4).ColumnWidth = 30.6667
$ws.Columns.Item(5).ColumnWidth = 27.6667

# Move the active selection to E1
$ws.Range("E1").Select()
